# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Force the cell to remain plain text (matches original inlineStr cells)
    # even when the value looks numeric (e.g. "1.00" or "5.36").
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "58.202.92"
Set-TextValue $ws.Range("E2") "  -1.13%  "
Set-TextValue $ws.Range("D3") "2.476.36"
Set-TextValue $ws.Range("E3") "  -1.80%  "
Set-TextValue $ws.Range("E4") "  -0.12%  "
Set-TextValue $ws.Range("D5") "520.48"
Set-TextValue $ws.Range("E5") "  -3.09%  "
Set-TextValue $ws.Range("D6") "131.78"
Set-TextValue $ws.Range("E6") "  -4.03%  "
Set-TextValue $ws.Range("E7") "  -0.01%  "
Set-TextValue $ws.Range("D8") "0.559"
Set-TextValue $ws.Range("E8") "  -1.39%  "
Set-TextValue $ws.Range("D9") "0.0996"
Set-TextValue $ws.Range("E9") "  -1.52%  "
Set-TextValue $ws.Range("E10") "  -0.77%  "
Set-TextValue $ws.Range("D11") "5.36"
Set-TextValue $ws.Range("E11") "  +0.25%  "
Set-TextValue $ws.Range("D12") "0.343"
Set-TextValue $ws.Range("E12") "  -1.20%  "
Set-TextValue $ws.Range("D13") "2.916.17"
Set-TextValue $ws.Range("E13") "  -1.71%  "
Set-TextValue $ws.Range("D14") "58.142.36"
Set-TextValue $ws.Range("E14") "  -1.26%  "
Set-TextValue $ws.Range("D15") "22.10"
Set-TextValue $ws.Range("E15") "  -4.23%  "
Set-TextValue $ws.Range("E16") "  -1.80%  "
Set-TextValue $ws.Range("D17") "2.481.07"
Set-TextValue $ws.Range("E17") "  -1.80%  "
Set-TextValue $ws.Range("D18") "10.84"
Set-TextValue $ws.Range("E18") "  -2.66%  "
Set-TextValue $ws.Range("D19") "4.19"
Set-TextValue $ws.Range("E19") "  -2.18%  "
Set-TextValue $ws.Range("E20") "  -0.81%  "
Set-TextValue $ws.Range("E21") "  -0.19%  "
Set-TextValue $ws.Range("D22") "5.77"
Set-TextValue $ws.Range("E22") "  -2.69%  "
Set-TextValue $ws.Range("D23") "64.14"
Set-TextValue $ws.Range("E23") "  -1.89%  "
Set-TextValue $ws.Range("E24") "  -2.54%  "
Set-TextValue $ws.Range("D25") "1.00"
Set-TextValue $ws.Range("E25") "  -0.10%  "
Set-TextValue $ws.Range("E26") "  -3.05%  "
Set-TextValue $ws.Range("D27") "7.41"
Set-TextValue $ws.Range("E27") "  -2.73%  "
Set-TextValue $ws.Range("D28") "0.0₃0755"
Set-TextValue $ws.Range("E28") "  -2.13%  "
Set-TextValue $ws.Range("E29") "  -4.00%  "
Set-TextValue $ws.Range("D30") "6.35"
Set-TextValue $ws.Range("E30") "  -6.09%  "
Set-TextValue $ws.Range("E31") "  +1.34%  "
Set-TextValue $ws.Range("D32") "166.35"
Set-TextValue $ws.Range("E32") "  +0.32%  "
Set-TextValue $ws.Range("E33") "  -0.01%  "
Set-TextValue $ws.Range("E34") "  +0.05%  "
Set-TextValue $ws.Range("D35") "18.17"
Set-TextValue $ws.Range("E35") "  -1.41%  "
Set-TextValue $ws.Range("E36") "  -10.47%  "
Set-TextValue $ws.Range("E37") "  -2.54%  "
Set-TextValue $ws.Range("E38") "  -3.67%  "
Set-TextValue $ws.Range("D39") "0.796"
Set-TextValue $ws.Range("E39") "  -2.45%  "
Set-TextValue $ws.Range("B40") "Filecoin"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D40") "3.49"
Set-TextValue $ws.Range("E40") "  -3.83%  "
Set-TextValue $ws.Range("B41") "Bittensor"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D41") "276.91"
Set-TextValue $ws.Range("E41") "  -3.14%  "
Set-TextValue $ws.Range("D42") "5.07"
Set-TextValue $ws.Range("E42") "  -1.83%  "
Set-TextValue $ws.Range("E43") "  -1.17%  "
Set-TextValue $ws.Range("D44") "126.52"
Set-TextValue $ws.Range("E44") "  -4.38%  "
Set-TextValue $ws.Range("D45") "0.0909"
Set-TextValue $ws.Range("E45") "  -2.16%  "
Set-TextValue $ws.Range("D46") "0.0492"
Set-TextValue $ws.Range("E46") "  -3.19%  "
Set-TextValue $ws.Range("E47") "  -2.68%  "
Set-TextValue $ws.Range("D48") "17.22"
Set-TextValue $ws.Range("E48") "  -0.43%  "
Set-TextValue $ws.Range("D49") "1.744.58"
Set-TextValue $ws.Range("E49") "  -1.37%  "
Set-TextValue $ws.Range("D50") "0.973"
Set-TextValue $ws.Range("E50") "  -2.02%  "
Set-TextValue $ws.Range("E51") "  -1.23%  "
